$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2.2922916666666668
$ws.Range("B5").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C5").Value = "House of Anubis Season 2 (Audiovisual, English, New):24; Bleach (Text with visuals, Japanese, Familiar):27;"
$ws.Range("D5").Value = "Watched children's shows from my childhood and their continuations and spin-offs with, and read simple manga."

$ws.Range("B5").Select()
